$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# --- Update the header texts in A1:C1 -------------------------------------
# A1 keeps the Kyrgyz title, B1 the Russian title, C1 the English title -
# all three get a small punctuation tweak (the "16.5.1.1a." prefix loses its
# trailing period and gets a space instead).
$ws.Range("A1").Value = '16.5.1.1a "Аткаруу бийлигинин мамлекеттик органдарындагы жана жергиликтүү өз алдынча башкаруу органдарындагы коррупциянын деңгээли жөнүндө жеке түшүнүк" индекси'
$ws.Range("B1").Value = '16.5.1.1a Индекс "Личное представление об уровне коррупции в государственных органах исполнительной власти и органах местного самоуправления'''''
$ws.Range("C1").Value = '16.5.1.1a Index "Personal views about the level of corruption in executive government authorities and local government'''''

# --- Add the new 2020 data column (column I) -------------------------------
# Column I mirrors the formatting of column H (the 2019 column): the header
# cell copies H4's style outright, while the data cells copy H's font/border
# but additionally get the "0.0" number format applied.

$ws.Range("H4").Copy() | Out-Null
$ws.Range("I4").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("I4").Value = 2020

$dataRows = @{
    5  = 12.3
    6  = 40.3
    7  = 36.2
    8  = 44.3
    9  = 36
    10 = 2.7
    11 = 32.9
    12 = 11.3
    13 = -18.2
    14 = 33
}

foreach ($row in $dataRows.Keys) {
    $ws.Range("H$row").Copy() | Out-Null
    $ws.Range("I$row").PasteSpecial($xlPasteFormats) | Out-Null
    $ws.Range("I$row").NumberFormat = "0.0"
    $ws.Range("I$row").Value = $dataRows[$row]
}

# --- Move the active selection ---------------------------------------------
$ws.Range("F16").Select() | Out-Null
